# Insert 5 new price rows (weekly update) before row 954 in the single
# worksheet, shifting all existing rows 954-1049 down to 959-1054.
# This matches the commit "Fruta / hortaliza, semanal" which adds a new
# week's worth of market observations for "Poroto verde".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 954
$numNewRows = 5

# Shift existing rows 954..end down by $numNewRows, carrying formatting
# (date style on column D, etc.) from the row being pushed down - this
# mirrors Excel's native "Insert Copied Cells" / "Shift cells down" UX.
$ws.Range("A" + $startRow + ":R" + ($startRow + $numNewRows - 1)).Insert()

# New row data, in column order A..R
$rows = @(
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44918, 13, 100112031, 'Poroto verde', 'Magnum', 'Primera', 400, 25000, 27000, 25850, '$/saco 25 kilos', 'Región Metropolitana', 1034, 25, 'Hortaliza'),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44918, 13, 100112031, 'Poroto verde', 'Magnum', 'Primera', 580, 25000, 27000, 25897, '$/saco 25 kilos', "Región de O'Higgins", 1036, 25, 'Hortaliza'),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44918, 13, 100112031, 'Poroto verde', 'Magnum', 'Segunda', 270, 15000, 15000, 15000, '$/saco 25 kilos', "Región de O'Higgins", 600, 25, 'Hortaliza'),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44918, 13, 100112031, 'Poroto verde', 'Sin especificar', 'Primera', 800, 30000, 35000, 32812, '$/malla 25 kilos', 'Región de Coquimbo', 1312, 25, 'Hortaliza'),
    @(6, 'Mercado Mayorista Lo Valledor de Santiago', 'Metropolitana', 44918, 13, 100112031, 'Poroto verde', 'Sin especificar', 'Segunda', 250, 25000, 25000, 25000, '$/malla 25 kilos', 'Región de Coquimbo', 1000, 25, 'Hortaliza')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
